# Changing cdr from uniform distribution to beta distribution
# Update the dependent constants on the "constants" sheet that flow
# from the change in distribution used for cdr_adjustment.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

$ws.Range("B2").Value  = 16.48368279065028   # tb_n_contact
$ws.Range("B3").Value  = 0.6596736217947906  # cdr_adjustment
$ws.Range("B5").Value  = 1899.68673035834    # start_time
$ws.Range("B6").Value  = 113795.6479117848   # susceptible_fully
$ws.Range("B14").Value = 0.7335900370457137  # tb_prop_casefatality_untreated_smearpos
$ws.Range("B15").Value = 1.102460887445562   # tb_multiplier_treated_protection
$ws.Range("B16").Value = 3.357761185341813   # tb_timeperiod_activeuntreated
$ws.Range("B17").Value = 0.8534673593383884  # active
